# Updates cryptos list values (Price / Volume(1h)) per commit
# 'Updated cryptos list on Tue Nov 14 17:55:50 UTC 2023 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $text) {
    # Force the cell to remain plain text so numeric-looking strings
    # (e.g. '244.15') are not silently converted into floating point
    # numbers (which would lose precision / change representation).
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "36.215.20"
Set-TextCell $ws "E2" "  -1.30%  "

Set-TextCell $ws "D3" "2.033.92"
Set-TextCell $ws "E3" "  -2.26%  "

Set-TextCell $ws "E4" "  +0.24%  "

Set-TextCell $ws "D5" "244.15"
Set-TextCell $ws "E5" "  -0.19%  "

Set-TextCell $ws "D6" "0.658"
Set-TextCell $ws "E6" "  +1.38%  "

Set-TextCell $ws "E7" "  +0.07%  "

Set-TextCell $ws "D8" "54.25"
Set-TextCell $ws "E8" "  +1.07%  "

Set-TextCell $ws "D9" "59.05"
Set-TextCell $ws "E9" "  +0.67%  "

Set-TextCell $ws "E10" "  -0.62%  "

Set-TextCell $ws "D11" "0.0735"
Set-TextCell $ws "E11" "  -3.23%  "

Set-TextCell $ws "E12" "  -4.05%  "

Set-TextCell $ws "D13" "0.894"
Set-TextCell $ws "E13" "  +1.19%  "

Set-TextCell $ws "D14" "14.15"
Set-TextCell $ws "E14" "  -5.04%  "

Set-TextCell $ws "D15" "2.336.91"
Set-TextCell $ws "E15" "  -2.07%  "

Set-TextCell $ws "D16" "5.31"
Set-TextCell $ws "E16" "  -2.94%  "

Set-TextCell $ws "D17" "2.033.30"
Set-TextCell $ws "E17" "  -2.49%  "

Set-TextCell $ws "D18" "17.36"
Set-TextCell $ws "E18" "  +0.98%  "

Set-TextCell $ws "D19" "36.120.17"
Set-TextCell $ws "E19" "  -1.51%  "

Set-TextCell $ws "D20" "71.19"
Set-TextCell $ws "E20" "  -1.63%  "

Set-TextCell $ws "D21" "0.0₃0850"
Set-TextCell $ws "E21" "  -2.81%  "

Set-TextCell $ws "D22" "235.57"
Set-TextCell $ws "E22" "  -1.85%  "

Set-TextCell $ws "D23" "5.15"
Set-TextCell $ws "E23" "  -4.89%  "

Set-TextCell $ws "E24" "  +0.12%  "

Set-TextCell $ws "E25" "  -2.10%  "

Set-TextCell $ws "D26" "2.27"
Set-TextCell $ws "E26" "  +5.73%  "

Set-TextCell $ws "D27" "9.17"
Set-TextCell $ws "E27" "  -6.54%  "

Set-TextCell $ws "D28" "163.03"
Set-TextCell $ws "E28" "  -2.44%  "

Set-TextCell $ws "D29" "19.80"
Set-TextCell $ws "E29" "  -4.06%  "

Set-TextCell $ws "E30" "  -2.16%  "

Set-TextCell $ws "E31" "  -0.83%  "

Set-TextCell $ws "D32" "4.92"
Set-TextCell $ws "E32" "  -6.87%  "

Set-TextCell $ws "D33" "0.0595"
Set-TextCell $ws "E33" "  -1.48%  "

Set-TextCell $ws "D34" "4.34"
Set-TextCell $ws "E34" "  -6.56%  "

Set-TextCell $ws "D35" "0.0893"
Set-TextCell $ws "E35" "  +8.81%  "

Set-TextCell $ws "E36" "  +0.17%  "

Set-TextCell $ws "D37" "1.82"
Set-TextCell $ws "E37" "  -1.47%  "

Set-TextCell $ws "D38" "2.19"
Set-TextCell $ws "E38" "  -8.02%  "

Set-TextCell $ws "E39" "  +4.01%  "

Set-TextCell $ws "E40" "  -5.10%  "

Set-TextCell $ws "E41" "  +1.79%  "

Set-TextCell $ws "E42" "  -2.38%  "

Set-TextCell $ws "D43" "1.09"
Set-TextCell $ws "E43" "  -4.77%  "

Set-TextCell $ws "D44" "0.0898"
Set-TextCell $ws "E44" "  -5.23%  "

Set-TextCell $ws "D45" "1.397.34"
Set-TextCell $ws "E45" "  +2.36%  "

Set-TextCell $ws "D46" "91.73"
Set-TextCell $ws "E46" "  -4.37%  "

Set-TextCell $ws "D47" "7.42"
Set-TextCell $ws "E47" "  +1.54%  "

Set-TextCell $ws "D48" "15.37"
Set-TextCell $ws "E48" "  -3.63%  "

Set-TextCell $ws "E49" "  +1.44%  "

Set-TextCell $ws "E50" "  -7.40%  "

Set-TextCell $ws "D51" "45.63"
Set-TextCell $ws "E51" "  +1.59%  "

